$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-align the existing data rows.
#    Rows 2-6 (bigger font, fontId=1) and rows 7-22 (default font) get
#    their id/name columns (A,B) centered. One stray cell, C8, picked up
#    a left alignment in the original edit as well.
# ---------------------------------------------------------------------
# -4131 = xlLeft, -4108 = xlCenter
$ws.Range("C8").HorizontalAlignment = -4131
$ws.Range("A7:B22").HorizontalAlignment = -4108
$ws.Range("A2:B6").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 2) Append the three new rows of data (social-media logo checks),
#    re-using the same "center" style that the rest of the id/name
#    columns in the lower block now have.
# ---------------------------------------------------------------------
$ws.Range("A23").Value = 800
$ws.Range("B23").Value = "check twitter logo"
$ws.Range("C23").Value = "while clicking on that logo open twitter page."

$ws.Range("A24").Value = 900
$ws.Range("B24").Value = "check fb logo"
$ws.Range("C24").Value = "while clicking on that logo open acebook page."

$ws.Range("A25").Value = 1000
$ws.Range("B25").Value = "check linkdin logo"
$ws.Range("C25").Value = "while clicking on that logo open linkdin page page."

$ws.Range("A23:B25").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 3) Three trailing blank rows, styled the same as column A above them.
# ---------------------------------------------------------------------
$ws.Range("A26").HorizontalAlignment = -4108
$ws.Range("A27").HorizontalAlignment = -4108
$ws.Range("A28").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 4) Selection moves to the freshly extended id column.
# ---------------------------------------------------------------------
$ws.Range("A2:A28").Select()
